$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.159.77'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '2.478.08'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''520.32'
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D6").Value = '''131.27'
$ws.Range("E6").Value = '  -3.92%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").Value = '''0.559'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = '''0.0994'
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = '''5.35'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '''0.344'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = '2.921.52'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").Value = '58.120.84'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '''22.39'
$ws.Range("E15").Value = '  -2.17%  '
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").Value = '2.481.06'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = '''10.86'
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("E19").Value = '  -1.79%  '
$ws.Range("D20").Value = '''319.70'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '''5.77'
$ws.Range("E22").Value = '  -3.01%  '
$ws.Range("D23").Value = '''64.19'
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '''7.35'
$ws.Range("E27").Value = '  -2.35%  '
$ws.Range("D28").Value = '0.0₃0756'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("E29").Value = '  -3.32%  '
$ws.Range("D30").Value = '''6.33'
$ws.Range("E30").Value = '  -5.30%  '
$ws.Range("D31").Value = '''165.99'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").Value = '''1.17'
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '''18.11'
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("D36").Value = '''1.32'
$ws.Range("E36").Value = '  -9.50%  '
$ws.Range("D37").Value = '''3.99'
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("E38").Value = '  -3.60%  '
$ws.Range("D39").Value = '''0.790'
$ws.Range("E39").Value = '  -2.53%  '
$ws.Range("D40").Value = '''3.48'
$ws.Range("E40").Value = '  -3.28%  '
$ws.Range("D41").Value = '''275.81'
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").Value = '''5.03'
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").Value = '''0.595'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").Value = '''127.39'
$ws.Range("E44").Value = '  -3.04%  '
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '''0.0489'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '''0.0214'
$ws.Range("E47").Value = '  -2.32%  '
$ws.Range("D48").Value = '''17.13'
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = '1.741.66'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  -1.01%  '
